$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = -13.6303
$ws.Range("B9").Value = 5.028999999999998
$ws.Range("C12").Value = -10.6802
$ws.Range("B13").Value = 6.486199999999997
$ws.Range("C14").Value = -12.85959999999999
$ws.Range("B16").Value = 6.157800000000003
$ws.Range("B18").Value = 7.325599999999997
$ws.Range("C19").Value = -11.94830000000001
$ws.Range("B20").Value = 8.797599999999999
$ws.Range("B26").Value = 5.623000000000006
$ws.Range("C26").Value = -13.4281
$ws.Range("B27").Value = 5.554000000000005
$ws.Range("C27").Value = -12.8301
$ws.Range("B29").Value = 4.834499999999997
$ws.Range("C29").Value = -10.68630000000001
$ws.Range("B35").Value = 8.341300000000006
$ws.Range("B36").Value = 8.722600000000009
$ws.Range("C37").Value = -13.3788
$ws.Range("C38").Value = -13.0529
$ws.Range("B45").Value = 5.317000000000002
$ws.Range("C47").Value = -12.0501
$ws.Range("C51").Value = -12.0366
$ws.Range("C52").Value = -11.4358
$ws.Range("B55").Value = 6.797199999999992
$ws.Range("C55").Value = -14.1788
$ws.Range("B57").Value = 5.140899999999997
$ws.Range("B69").Value = 5.927199999999995
$ws.Range("C69").Value = -11.4875
$ws.Range("C70").Value = -11.4286
$ws.Range("B76").Value = 5.471000000000004
$ws.Range("C76").Value = -11.27890000000001
$ws.Range("B78").Value = 9.829800000000004
$ws.Range("C81").Value = -12.8507
$ws.Range("B82").Value = 5.309399999999999
$ws.Range("B83").Value = 5.402199999999997
$ws.Range("C83").Value = -13.72779999999999
$ws.Range("B93").Value = 5.552899999999999
$ws.Range("C94").Value = -10.2002
$ws.Range("B97").Value = 5.986599999999999
$ws.Range("C100").Value = -12.50929999999999
$ws.Range("C102").Value = -13.249
